$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: downgrade row 96 (previously the "latest" / highlighted row) from
# the "Good" highlight style to the normal "Neutral" data-row style, matching
# every other historical row (81-95). Values/formulas stay the same for
# A96:H96; only the style (and the trailing helper cells I96/J96) change.

$ws.Range("A96").Style = "Neutral"
$ws.Range("A96").NumberFormat = "d-mmm"

$ws.Range("B96").Style = "Neutral"
$ws.Range("B96").NumberFormat = "General"

$ws.Range("C96").Style = "Neutral"
$ws.Range("C96").NumberFormat = "0"

$ws.Range("D96").Style = "Neutral"
$ws.Range("D96").NumberFormat = "General"

$ws.Range("E96").Style = "Neutral"
$ws.Range("E96").NumberFormat = "General"

$ws.Range("G96").Style = "Neutral"
$ws.Range("G96").NumberFormat = "General"

$ws.Range("H96").Style = "Neutral"
$ws.Range("H96").NumberFormat = "m/d/yyyy"

# I96 no longer carries the "days left" helper formula - blank it out but
# keep it formatted like the rest of the Neutral rows (e.g. I95).
$ws.Range("I96").ClearContents()
$ws.Range("I96").Style = "Neutral"
$ws.Range("I96").NumberFormat = "0"

# J96 no longer exists at all once row 96 stops being the highlighted row.
$ws.Range("J96").Clear()

# --- Step 2: add the new day of data as row 97, taking over the "Good"
# highlight style and the helper columns I/J that row 96 used to have.

$ws.Range("A97").Style = "Good"
$ws.Range("A97").NumberFormat = "d-mmm"
$ws.Range("A97").Value = 44282

$ws.Range("B97").Style = "Good"
$ws.Range("B97").NumberFormat = "General"
$ws.Range("B97").Value = 2951

$ws.Range("C97").Style = "Good"
$ws.Range("C97").NumberFormat = "0"
$ws.Range("C97").Formula = "=(AVERAGE(B91:B97))"

$ws.Range("D97").Style = "Good"
$ws.Range("D97").NumberFormat = "General"
$ws.Range("D97").Formula = "=AVERAGE(B84:B97)"

$ws.Range("E97").Style = "Good"
$ws.Range("E97").NumberFormat = "General"
$ws.Range("E97").Formula = "=(E96-B97)"

$ws.Range("G97").Style = "Good"
$ws.Range("G97").NumberFormat = "General"
$ws.Range("G97").Formula = "=E97/C97"

$ws.Range("H97").Style = "Good"
$ws.Range("H97").NumberFormat = "m/d/yyyy"
$ws.Range("H97").Formula = "=A97+G97"

$ws.Range("I97").Style = "Good"
$ws.Range("I97").NumberFormat = "0"
$ws.Range("I97").Formula = "=E97/84"

$ws.Range("J97").Style = "Good"
$ws.Range("J97").NumberFormat = "General"
$ws.Range("J97").Value = "daily rate to achieve June 20 target"

# --- Step 3: update the view so the selection matches the saved workbook
# state (scrolled further down, new active cell).
$ws.Range("G110").Select()
